# Crear segunda prueba de escritorio
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Rename the original sheet, matching the target names.
$ws1.Name = "Ejercicio 1"

# Add a new worksheet right after "Ejercicio 1" for the second desktop test.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "Ejercicio 2"

# --- Apply cell formatting by copying it from the already-styled cells on
# --- "Ejercicio 1" (keeps the shared style/border/font table unchanged
# --- instead of growing new duplicate style entries).

# Style used by the header rows + the "Proceso/SubProceso" column (bold, bordered).
$boldBordered = $ws1.Range("B3")
$boldBordered.Copy()
$newSheet.Range("C2:G2").PasteSpecial(-4122)
$newSheet.Range("C3:C12").PasteSpecial(-4122)
$newSheet.Range("C16:G16").PasteSpecial(-4122)
$newSheet.Range("C17:C30").PasteSpecial(-4122)

# Style used by the plain data cells (bordered, not bold).
$plainBordered = $ws1.Range("C4")
$plainBordered.Copy()
$newSheet.Range("D3:G12").PasteSpecial(-4122)
$newSheet.Range("H16").PasteSpecial(-4122)
$newSheet.Range("D17:H30").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Fill in the values for both trace tables.
$newSheet.Range("C2").Value = "Proceso/SubProcesoLinea(inst)"
$newSheet.Range("D2").Value = "division"
$newSheet.Range("E2").Value = "modulo"
$newSheet.Range("F2").Value = "num1"
$newSheet.Range("G2").Value = "num2"
$newSheet.Range("C3").Value = "1:EJERCICIO_26(1)"
$newSheet.Range("D3").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E3").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F3").Value = "<<Variable no inicializada (NUM1).>>"
$newSheet.Range("G3").Value = "<<Variable no inicializada (NUM2).>>"
$newSheet.Range("C4").Value = "1:EJERCICIO_27(1)"
$newSheet.Range("D4").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E4").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F4").Value = "<<Variable no inicializada (NUM1).>>"
$newSheet.Range("G4").Value = "<<Variable no inicializada (NUM2).>>"
$newSheet.Range("C5").Value = "1:EJERCICIO_29(1)"
$newSheet.Range("D5").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E5").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F5").Value = "<<Variable no inicializada (NUM1).>>"
$newSheet.Range("G5").Value = "<<Variable no inicializada (NUM2).>>"
$newSheet.Range("C6").Value = "1:EJERCICIO_210(1)"
$newSheet.Range("D6").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E6").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F6").Value = "<<Variable no inicializada (NUM1).>>"
$newSheet.Range("G6").Value = "<<Variable no inicializada (NUM2).>>"
$newSheet.Range("C7").Value = "1:EJERCICIO_211(1)"
$newSheet.Range("D7").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E7").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F7").Value = 2
$newSheet.Range("G7").Value = "<<Variable no inicializada (NUM2).>>"
$newSheet.Range("C8").Value = "1:EJERCICIO_212(1)"
$newSheet.Range("D8").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E8").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F8").Value = 2
$newSheet.Range("G8").Value = "<<Variable no inicializada (NUM2).>>"
$newSheet.Range("C9").Value = "1:EJERCICIO_214(1)"
$newSheet.Range("D9").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E9").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F9").Value = 2
$newSheet.Range("G9").Value = 0
$newSheet.Range("C10").Value = "1:EJERCICIO_215(1)"
$newSheet.Range("D10").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E10").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F10").Value = 2
$newSheet.Range("G10").Value = 0
$newSheet.Range("C11").Value = "1:EJERCICIO_222(1)"
$newSheet.Range("D11").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E11").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F11").Value = 2
$newSheet.Range("G11").Value = 0
$newSheet.Range("C12").Value = "1:EJERCICIO_223(1)"
$newSheet.Range("D12").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E12").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F12").Value = 2
$newSheet.Range("G12").Value = 0
$newSheet.Range("C16").Value = "Proceso/SubProcesoLinea(inst)"
$newSheet.Range("D16").Value = "division"
$newSheet.Range("E16").Value = "modulo"
$newSheet.Range("F16").Value = "num1"
$newSheet.Range("G16").Value = "num2"
$newSheet.Range("C17").Value = "1:EJERCICIO_26(1)"
$newSheet.Range("D17").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E17").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F17").Value = "<<Variable no inicializada (NUM1).>>"
$newSheet.Range("G17").Value = "<<Variable no inicializada (NUM2).>>"
$newSheet.Range("C18").Value = "1:EJERCICIO_27(1)"
$newSheet.Range("D18").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E18").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F18").Value = "<<Variable no inicializada (NUM1).>>"
$newSheet.Range("G18").Value = "<<Variable no inicializada (NUM2).>>"
$newSheet.Range("C19").Value = "1:EJERCICIO_29(1)"
$newSheet.Range("D19").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E19").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F19").Value = "<<Variable no inicializada (NUM1).>>"
$newSheet.Range("G19").Value = "<<Variable no inicializada (NUM2).>>"
$newSheet.Range("C20").Value = "1:EJERCICIO_210(1)"
$newSheet.Range("D20").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E20").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F20").Value = "<<Variable no inicializada (NUM1).>>"
$newSheet.Range("G20").Value = "<<Variable no inicializada (NUM2).>>"
$newSheet.Range("C21").Value = "1:EJERCICIO_211(1)"
$newSheet.Range("D21").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E21").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F21").Value = 10
$newSheet.Range("G21").Value = "<<Variable no inicializada (NUM2).>>"
$newSheet.Range("C22").Value = "1:EJERCICIO_212(1)"
$newSheet.Range("D22").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E22").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F22").Value = 10
$newSheet.Range("G22").Value = "<<Variable no inicializada (NUM2).>>"
$newSheet.Range("C23").Value = "1:EJERCICIO_214(1)"
$newSheet.Range("D23").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E23").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F23").Value = 10
$newSheet.Range("G23").Value = 2
$newSheet.Range("C24").Value = "1:EJERCICIO_216(1)"
$newSheet.Range("D24").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E24").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F24").Value = 10
$newSheet.Range("G24").Value = 2
$newSheet.Range("C25").Value = "1:EJERCICIO_217(1)"
$newSheet.Range("D25").Value = "<<Variable no inicializada (DIVISION).>>"
$newSheet.Range("E25").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F25").Value = 10
$newSheet.Range("G25").Value = 2
$newSheet.Range("C26").Value = "1:EJERCICIO_218(1)"
$newSheet.Range("D26").Value = 5
$newSheet.Range("E26").Value = "<<Variable no inicializada (MODULO).>>"
$newSheet.Range("F26").Value = 10
$newSheet.Range("G26").Value = 2
$newSheet.Range("C27").Value = "1:EJERCICIO_220(1)"
$newSheet.Range("D27").Value = 5
$newSheet.Range("E27").Value = 0
$newSheet.Range("F27").Value = 10
$newSheet.Range("G27").Value = 2
$newSheet.Range("C28").Value = "1:EJERCICIO_221(1)"
$newSheet.Range("D28").Value = 5
$newSheet.Range("E28").Value = 0
$newSheet.Range("F28").Value = 10
$newSheet.Range("G28").Value = 2
$newSheet.Range("C29").Value = "1:EJERCICIO_222(1)"
$newSheet.Range("D29").Value = 5
$newSheet.Range("E29").Value = 0
$newSheet.Range("F29").Value = 10
$newSheet.Range("G29").Value = 2
$newSheet.Range("C30").Value = "1:EJERCICIO_223(1)"
$newSheet.Range("D30").Value = 5
$newSheet.Range("E30").Value = 0
$newSheet.Range("F30").Value = 10
$newSheet.Range("G30").Value = 2

# --- Column widths roughly matching the source workbook's autosized columns.
$newSheet.Columns.Item(3).AutoFit()
$newSheet.Columns.Item(4).AutoFit()
$newSheet.Columns.Item(5).AutoFit()
$newSheet.Columns.Item(6).AutoFit()
$newSheet.Columns.Item(7).AutoFit()

# --- Match the saved selection/active sheet from the target workbook.
$newSheet.Activate()
$newSheet.Range("D6").Select()
